$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for 178469f5-... row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-24 16:47:13"

# zh-cn sheet: 178469f5-... row - Correspond Handoff Datetime (H3) and Correspond Handback DateTime (K3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-24 16:47:04"
$wsZhCn.Range("K3").Value = "2016-08-24 16:47:32"

# de-de sheet: 178469f5-... row - Correspond Handback DateTime (K3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K3").Value = "2016-08-24 16:47:39"
